# Update betting-odds figures per 2024-11-05 FlashScore refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("G2").Value = 2.1
$ws.Range("H2").Value = 2.88
$ws.Range("I2").Value = 4.33
$ws.Range("J2").Value = 3
$ws.Range("L2").Value = 5.5
$ws.Range("M2").Value = 1.17
$ws.Range("N2").Value = 4.75
$ws.Range("S2").Value = 1.75
$ws.Range("T2").Value = 2.05
$ws.Range("W2").Value = 4.75
$ws.Range("X2").Value = 8
$ws.Range("Y2").Value = 11
$ws.Range("Z2").Value = 19
$ws.Range("AA2").Value = 26
$ws.Range("AC2").Value = 4.75
$ws.Range("AE2").Value = 26
$ws.Range("AG2").Value = 7.5
$ws.Range("AH2").Value = 19
$ws.Range("AI2").Value = 17
$ws.Range("AJ2").Value = 51
$ws.Range("AK2").Value = 51
$ws.Range("AN2").Value = 3.75
$ws.Range("AO2").Value = 13
$ws.Range("AP2").Value = 34
$ws.Range("AR2").Value = 101
$ws.Range("AS2").Value = 451
$ws.Range("AW2").Value = 5.5
$ws.Range("AX2").Value = 29
$ws.Range("AZ2").Value = 126
$ws.Range("BA2").Value = 201

# Row 3
$ws.Range("G3").Value = 2.9
$ws.Range("H3").Value = 2.55
$ws.Range("I3").Value = 3.1
$ws.Range("M3").Value = 1.18
$ws.Range("N3").Value = 4.5
$ws.Range("O3").Value = 1.8
$ws.Range("P3").Value = 1.91
$ws.Range("Q3").Value = 3.6
$ws.Range("R3").Value = 1.29
$ws.Range("S3").Value = 1.8
$ws.Range("T3").Value = 2
$ws.Range("Y3").Value = 13
$ws.Range("AA3").Value = 34
$ws.Range("AE3").Value = 23
$ws.Range("AG3").Value = 6
$ws.Range("AH3").Value = 13
$ws.Range("AN3").Value = 4.33
$ws.Range("AT3").Value = 1.91

# Row 4
$ws.Range("S4").Value = 1.5

# Row 5
$ws.Range("Q5").Value = 1.7
$ws.Range("R5").Value = 2.1

# Row 6
$ws.Range("M6").Value = 1.03
$ws.Range("N6").Value = 15
$ws.Range("S6").Value = 1.3

# Row 7
$ws.Range("M7").Value = 1.13
$ws.Range("N7").Value = 6
$ws.Range("S7").Value = 1.62

# Row 8
$ws.Range("G8").Value = 1.75
$ws.Range("H8").Value = 3.4
$ws.Range("I8").Value = 5
$ws.Range("J8").Value = 2.5
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 5.5
$ws.Range("U8").Value = 2.2
$ws.Range("V8").Value = 1.62
$ws.Range("X8").Value = 7
$ws.Range("Z8").Value = 13
$ws.Range("AB8").Value = 34
$ws.Range("AC8").Value = 7
$ws.Range("AG8").Value = 10
$ws.Range("AH8").Value = 23
$ws.Range("AI8").Value = 17
$ws.Range("AJ8").Value = 51
$ws.Range("AK8").Value = 41
$ws.Range("AN8").Value = 3.5
$ws.Range("AO8").Value = 9.5
$ws.Range("AP8").Value = 26
$ws.Range("AQ8").Value = 34
$ws.Range("AR8").Value = 67
$ws.Range("AW8").Value = 6.5
$ws.Range("AX8").Value = 29

# Row 9
$ws.Range("O9").Value = 1.5
$ws.Range("P9").Value = 2.5

# Row 10
$ws.Range("G10").Value = 1.65
$ws.Range("H10").Value = 4.2
$ws.Range("J10").Value = 2.2
$ws.Range("S10").Value = 1.29
$ws.Range("T10").Value = 3.5
$ws.Range("AH10").Value = 26
$ws.Range("AN10").Value = 3.75
$ws.Range("AQ10").Value = 23
$ws.Range("AT10").Value = 3.5
$ws.Range("AZ10").Value = 81

# Row 11
$ws.Range("N11").Value = 13
